$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.478.24"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.44"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.89"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "56.56"
$ws.Range("E9").Value = "  +8.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.355"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0750"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0982"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +11.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.789"
$ws.Range("E14").Value = "  +8.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.165.84"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.01"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.883.13"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.475.64"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.23"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.12"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.99"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.16"
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.83"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.31"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.40"
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0606"
$ws.Range("E32").Value = "  +4.25%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("B34").Value = "BinanceUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  +19.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -16.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.850"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0747"
$ws.Range("E38").Value = "  +9.06%  "
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0228"
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.76"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.94"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.43"
$ws.Range("E43").Value = "  +18.80%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.08"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.311.38"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0810"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.36"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.51"
$ws.Range("E51").Value = "  -2.07%  "
